$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($l in $wb.LinkSources()) {
    $wb.BreakLink($l, 1)
}

$ws.Rows("1:1").Insert()

$ws.Range("C8").Select()
